$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow edits, then restore protection afterward.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclaimer text (cell A80).
[void]$ws.Range("A80").Replace("2021-03-18", "2021-03-19")

# Setting the multi-line cell value can trigger an automatic row-height
# resize; restore the row to its default auto-fit height.
$ws.Rows.Item(80).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-77.
$ws.Range("D2").Value = 0.06386155911810711
$ws.Range("E2").Value = -0.004480212395254313
$ws.Range("D3").Value = 0.03850437150466184
$ws.Range("E3").Value = 0.0155119402640036
$ws.Range("D4").Value = 0.03031665520694398
$ws.Range("E4").Value = -0.001603675450762898
$ws.Range("D5").Value = 0.02993008428242253
$ws.Range("E5").Value = -0.06235837940723277
$ws.Range("D6").Value = 0.02806581208918759
$ws.Range("E6").Value = -0.01592134475103091
$ws.Range("D7").Value = 0.02570366028198018
$ws.Range("E7").Value = 0.002780333837949289
$ws.Range("D8").Value = 0.02584709844081577
$ws.Range("E8").Value = -0.002679628591013894
$ws.Range("D9").Value = 0.1715344896819791
$ws.Range("E9").Value = -0.03154305200341001
$ws.Range("D10").Value = 0.0224735041554361
$ws.Range("E10").Value = 0.001923816852635651
$ws.Range("D11").Value = 0.02291195696719596
$ws.Range("E11").Value = -0.005361305361305302
$ws.Range("D12").Value = 0.02119052951251775
$ws.Range("E12").Value = -0.005928853754940788
$ws.Range("D13").Value = 0.02079700709330095
$ws.Range("E13").Value = -0.01052901900359526
$ws.Range("D14").Value = 0.01697605869197836
$ws.Range("E14").Value = 0.009707865168539387
$ws.Range("D15").Value = 0.01566968633520739
$ws.Range("E15").Value = 0.01138281757195414
$ws.Range("D16").Value = 0.01716849641098355
$ws.Range("E16").Value = 0.003688524590163889
$ws.Range("D17").Value = 0.01543011409119475
$ws.Range("E17").Value = 0.01330666871779118
$ws.Range("D18").Value = 0.01460568377519231
$ws.Range("E18").Value = 0.01486455896778072
$ws.Range("D19").Value = 0.01324408094532359
$ws.Range("E19").Value = -0.005632811124801962
$ws.Range("D20").Value = 0.01187599287955007
$ws.Range("E20").Value = -0.003372843789149038
$ws.Range("D21").Value = 0.01299090242205971
$ws.Range("E21").Value = 0.04123896346278078
$ws.Range("D22").Value = 0.01244487099117315
$ws.Range("E22").Value = -0.02858310626702998
$ws.Range("D23").Value = 0.01169478775868064
$ws.Range("E23").Value = 0.01702766179540705
$ws.Range("D24").Value = 0.01297479530020464
$ws.Range("E24").Value = -0.008330556481172979
$ws.Range("D25").Value = 0.01227701782657831
$ws.Range("E25").Value = 0.009750034525617846
$ws.Range("D26").Value = 0.01066020188963243
$ws.Range("E26").Value = 0.007682030728122902
$ws.Range("D27").Value = 0.010601877153652
$ws.Range("E27").Value = 0.0004797697105389798
$ws.Range("D28").Value = 0.01053437559693047
$ws.Range("E28").Value = 0.0005492349941151886
$ws.Range("D29").Value = 0.01041995145101537
$ws.Range("E29").Value = -0.007338464292106495
$ws.Range("D30").Value = 0.009141554621676875
$ws.Range("E30").Value = 0.0058237661590963
$ws.Range("D31").Value = 0.01052897123367647
$ws.Range("E31").Value = -0.002898550724637627
$ws.Range("D32").Value = 0.01080211411060807
$ws.Range("E32").Value = 0.02097753920045209
$ws.Range("D33").Value = 0.009447038903806933
$ws.Range("E33").Value = -0.001660123387548951
$ws.Range("D34").Value = 0.009785966657367631
$ws.Range("E34").Value = -0.001243118451429548
$ws.Range("D35").Value = 0.008305679771724783
$ws.Range("E35").Value = 0.002618041521220027
$ws.Range("D36").Value = 0.009217130932907332
$ws.Range("E36").Value = 0.004745896776744862
$ws.Range("D37").Value = 0.009399353345683399
$ws.Range("E37").Value = -0.01379932356257041
$ws.Range("D38").Value = 0.008708018720799087
$ws.Range("E38").Value = -0.027015186915888
$ws.Range("D39").Value = 0.007869600641080405
$ws.Range("E39").Value = 0.0222126467736723
$ws.Range("D40").Value = 0.008010156472847198
$ws.Range("E40").Value = 0.01811870290407258
$ws.Range("D41").Value = 0.00769962811834674
$ws.Range("E41").Value = 0.009788053949903786
$ws.Range("D42").Value = 0.008495998132591156
$ws.Range("E42").Value = -0.0396731158762309
$ws.Range("D43").Value = 0.008515665776119442
$ws.Range("E43").Value = 0.00310599197618755
$ws.Range("D44").Value = 0.007920973882365493
$ws.Range("E44").Value = -0.007192088702427157
$ws.Range("D45").Value = 0.007944032498915897
$ws.Range("E45").Value = -0.007107184018440349
$ws.Range("D46").Value = 0.008187419587578418
$ws.Range("E46").Value = -0.006709533128319789
$ws.Range("D47").Value = 0.007936911455569448
$ws.Range("E47").Value = -0.01794411689310427
$ws.Range("D48").Value = 0.006237016679581682
$ws.Range("E48").Value = 0.001957266351329201
$ws.Range("D49").Value = 0.007300086722015682
$ws.Range("E49").Value = -0.01415598290598286
$ws.Range("D50").Value = 0.006842559687006409
$ws.Range("E50").Value = -0.01189370005575174
$ws.Range("D51").Value = 0.006844509496494127
$ws.Range("E51").Value = -0.02336570140454308
$ws.Range("D52").Value = 0.00662969135554295
$ws.Range("E52").Value = -0.01745435016111707
$ws.Range("D53").Value = 0.006253801996041167
$ws.Range("E53").Value = -0.003660024400162776
$ws.Range("D54").Value = 0.006135372263243694
$ws.Range("E54").Value = 0.008166028767634392
$ws.Range("D55").Value = 0.005708915018549577
$ws.Range("E55").Value = 0.002303151798641156
$ws.Range("D56").Value = 0.005432338781433075
$ws.Range("E56").Value = -0.005680399500624156
$ws.Range("D57").Value = 0.006040679341601049
$ws.Range("E57").Value = -0.003368137420006856
$ws.Range("D58").Value = 0.005633762578946887
$ws.Range("E58").Value = -0.01124052004333698
$ws.Range("D59").Value = 0.005040257525750678
$ws.Range("E59").Value = -0.01230342275670671
$ws.Range("D60").Value = 0.005256855927538467
$ws.Range("E60").Value = 0.01149008224479897
$ws.Range("D61").Value = 0.004777033244908767
$ws.Range("E61").Value = -0.007985803016858917
$ws.Range("D62").Value = 0.004897921433147275
$ws.Range("E62").Value = 0.001107726391581298
$ws.Range("D63").Value = 0.004489648281284265
$ws.Range("E63").Value = -0.01570996978851968
$ws.Range("D64").Value = 0.004227186969372329
$ws.Range("E64").Value = 0.005936146317985003
$ws.Range("D65").Value = 0.004055603734453157
$ws.Range("E65").Value = -0.0007525083612038852
$ws.Range("D66").Value = 0.003736513173071851
$ws.Range("E66").Value = -0.006398039749523687
$ws.Range("D67").Value = 0.003913013318873095
$ws.Range("E67").Value = -0.002339789418952409
$ws.Range("D68").Value = 0.003031657043261843
$ws.Range("E68").Value = 0.0120380856507698
$ws.Range("D69").Value = 0.003350832378968702
$ws.Range("E69").Value = 0.06098440286896145
$ws.Range("D70").Value = 0.003355198256734679
$ws.Range("E70").Value = -0.009449694274596854
$ws.Range("D71").Value = 0.002663778857524814
$ws.Range("E71").Value = 0.01298453312965431
$ws.Range("D72").Value = 0.002157082713694828
$ws.Range("E72").Value = 0.009687561406956258
$ws.Range("D73").Value = 0.002139661589793697
$ws.Range("E73").Value = 0.01463975118366045
$ws.Range("D74").Value = 0.001861559414817242
$ws.Range("E74").Value = 0.005282572066123148
$ws.Range("D75").Value = 0.001688026370410353
$ws.Range("E75").Value = -0.004770992366412319
$ws.Range("D76").Value = 0.001677599128367339
$ws.Range("E76").Value = 0.02440749911567042
$ws.Range("E77").Value = -0.007365714766142339

# Restore sheet protection with the original password.
$ws.Protect("D382")

